$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Update the "总计" (Total) summary sheet: insert a new first data
#    row for "2022-Q1" (4 funds, 0.2 billion yuan held), pushing all
#    the existing quarter rows down by one.
# ------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")

for ($r = 6; $r -ge 2; $r--) {
    $nr = $r + 1
    $bVal = $totalWs.Range("B$r").Value()
    $cVal = $totalWs.Range("C$r").Value()
    $dVal = $totalWs.Range("D$r").Value()
    $totalWs.Range("A$nr").Value = $r - 1
    $totalWs.Range("B$nr").Value = $bVal
    $totalWs.Range("C$nr").Value = $cVal
    $totalWs.Range("D$nr").Value = $dVal
}

# Row 7 is brand new (didn't exist before) - copy the index-column
# formatting from row 6 so it keeps the same bold/bordered style.
$totalWs.Range("A6").Copy()
$totalWs.Range("A7").PasteSpecial(-4122)

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 4
$totalWs.Range("D2").Value = 0.2

# ------------------------------------------------------------------
# 2) Insert a brand new "2022-Q1" sheet right before "总计", holding
#    the fund-level breakdown for the quarter (same layout used by
#    the other quarterly sheets).
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($totalWs)
$q1.Name = "2022-Q1"

# Pull the header / index-column formatting (bold text + thin border,
# style shared with every other quarter sheet) from "2021-Q4" so the
# new sheet matches the existing look instead of using plain defaults.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$q1.Range("A2:A5").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# B:G hold text-formatted values in the source workbook (fund codes
# with leading zeros, percentages kept as fixed-decimal strings) -
# force text format before assigning so "011243" / "3.40" / "0.1340"
# round-trip exactly instead of being coerced to numbers.
$q1.Range("B2:G5").NumberFormat = "@"

$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "920002"
$q1.Range("C2").Value = "中金精选股票A"
$q1.Range("D2").Value = "3.40"
$q1.Range("E2").Value = "86.43"
$q1.Range("F2").Value = "3.94"
$q1.Range("G2").Value = "0.1340"
$q1.Range("H2").Value = 9

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "011243"
$q1.Range("C3").Value = "万家惠裕回报6个月持有期混合型证券投资基金A"
$q1.Range("D3").Value = "4.93"
$q1.Range("E3").Value = "23.04"
$q1.Range("F3").Value = "1.13"
$q1.Range("G3").Value = "0.0557"
$q1.Range("H3").Value = 4

$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "920922"
$q1.Range("C4").Value = "中金精选股票C"
$q1.Range("D4").Value = "0.14"
$q1.Range("E4").Value = "86.43"
$q1.Range("F4").Value = "3.94"
$q1.Range("G4").Value = "0.0055"
$q1.Range("H4").Value = 9

$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "011244"
$q1.Range("C5").Value = "万家惠裕回报6个月持有期混合型证券投资基金C"
$q1.Range("D5").Value = "0.14"
$q1.Range("E5").Value = "23.04"
$q1.Range("F5").Value = "1.13"
$q1.Range("G5").Value = "0.0016"
$q1.Range("H5").Value = 4

# Restore the originally active sheet.
$wb.Worksheets.Item("2020-Q4").Activate()
